$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final ticker list for column A (rows 2-387; row 1 is header "Sembol")
$values = @(
    'Sembol',
    'GOOGL',
    'GOOG',
    'MSFT',
    'AMZN',
    'META',
    'TSLA',
    'V',
    'ORCL',
    'MA',
    'AMD',
    'PLTR',
    'NFLX',
    'JPM-PD',
    'JPM-PC',
    'UNH',
    'BAC-PK',
    'BAC-PL',
    'NVO',
    'SAP',
    'IBM',
    'BAC-PE',
    'BML-PL',
    'BAC-PB',
    'AXP',
    'TMO',
    'CRM',
    'WFC-PY',
    'DIS',
    'WFC-PL',
    'APH',
    'ISRG',
    'ABT',
    'BX',
    'APP',
    'SHOP',
    'ACN',
    'BLK',
    'UBER',
    'DHR',
    'BKNG',
    'QCOM',
    'SPGI',
    'INTU',
    'UBS',
    'BBVA',
    'NOW',
    'COF',
    'BSX',
    'SONY',
    'SYK',
    'PANW',
    'ADBE',
    'WFC-PC',
    'PGR',
    'CRWD',
    'MELI',
    'KKR',
    'CVS',
    'SPOT',
    'ADP',
    'CEG',
    'CVNA',
    'LYG',
    'SNPS',
    'HOOD',
    'NKE',
    'MCO',
    'MRSH',
    'DASH',
    'CDNS',
    'B',
    'NTES',
    'ELV',
    'ORLY',
    'MS-PK',
    'AMT',
    'MS-PI',
    'ABNB',
    'TDG',
    'DB',
    'MS-PF',
    'MS-PE',
    'GM',
    'INFY',
    'SE',
    'USB-PP',
    'AON',
    'SNOW',
    'RELX',
    'NWG',
    'WBD',
    'TEL',
    'NET',
    'AJG',
    'AZO',
    'DUK-PA',
    'RACE',
    'RKT',
    'CTA-PB',
    'NXPI',
    'ADSK',
    'COIN',
    'NDAQ',
    'SRE',
    'IDXX',
    'MPLX',
    'TRI',
    'BIDU',
    'PYPL',
    'VST',
    'F',
    'ARGX',
    'MET',
    'EA',
    'SCHW-PD',
    'WDAY',
    'EW',
    'ARES',
    'AXON',
    'DDOG',
    'ALNY',
    'ROK',
    'CTA-PA',
    'HEI',
    'MSTR',
    'MSCI',
    'TTWO',
    'SPG-PJ',
    'EXC',
    'ROP',
    'JD',
    'EBAY',
    'MET-PE',
    'RKLB',
    'MET-PA',
    'EL',
    'CTSH',
    'TCOM',
    'LVS',
    'IQV',
    'RDDT',
    'MCHP',
    'CPRT',
    'XYZ',
    'HEI-A',
    'ALC',
    'ASTS',
    'A',
    'PRU',
    'PAYX',
    'CCI',
    'MDLN',
    'FICO',
    'VEEV',
    'GEHC',
    'FISV',
    'TEAM',
    'CPNG',
    'SATS',
    'ZS',
    'INSM',
    'NTRA',
    'MDB',
    'EXPE',
    'CHT',
    'ESLT',
    'PSA-PK',
    'ALL-PH',
    'SOFI',
    'ALL-PB',
    'FOXA',
    'HUM',
    'WTW',
    'FIS',
    'FOX',
    'VRSK',
    'FLUT',
    'MTD',
    'SYF',
    'DXCM',
    'TME',
    'STLA',
    'CSGP',
    'WIT',
    'BRO',
    'PPL',
    'EXE',
    'EFX',
    'ES',
    'FSLR',
    'STE',
    'AER',
    'DLTR',
    'AWK',
    'OMC',
    'AVB',
    'VLTO',
    'DLR-PK',
    'BR',
    'SQM',
    'AXIA-PC',
    'ILMN',
    'VRSN',
    'TPG',
    'TROW',
    'WAT',
    'CRDO',
    'LULU',
    'OWL',
    'CNC',
    'FUTU',
    'AFRM',
    'DLR-PJ',
    'TLK',
    'CYBR',
    'FWONK',
    'ALB',
    'PSLV',
    'FWONA',
    'CG',
    'RL',
    'UTHR',
    'GPN',
    'SSNC',
    'GMAB',
    'TWLO',
    'Q',
    'HL',
    'SBAC',
    'RCI',
    'CHKP',
    'PTC',
    'TOST',
    'GIB',
    'RIVN',
    'PODD',
    'TYL',
    'RVMD',
    'KTOS',
    'MRNA',
    'HIG-PG',
    'GRAB',
    'IOT',
    'U',
    'DKS',
    'HPQ',
    'CRCL',
    'FITBI',
    'XPEV',
    'IT',
    'PSNYW',
    'ALLY',
    'PNR',
    'PINS',
    'HUBS',
    'NWS',
    'IREN',
    'ZG',
    'TRMB',
    'MEDP',
    'Z',
    'TRU',
    'TTD',
    'INVH',
    'NLY',
    'MAA',
    'GFL',
    'WMG',
    'KSPI',
    'ROKU',
    'GEN',
    'NWSA',
    'H',
    'IONQ',
    'ONON',
    'DKNG',
    'AVAV',
    'KEY-PI',
    'EMA',
    'BBIO',
    'GH',
    'MLI',
    'HMY',
    'PFGC',
    'ULS',
    'GDDY',
    'ERIE',
    'ARCC',
    'ICLR',
    'FIG',
    'RGC',
    'KRMN',
    'W',
    'AKAM',
    'PEN',
    'CELH',
    'DPZ',
    'BBY',
    'EMBJ',
    'SBSW',
    'LOGI',
    'GWRE',
    'EG',
    'SOLV',
    'BILI',
    'RBRK',
    'FIGR',
    'PAA',
    'NLY-PG',
    'AMH',
    'RVTY',
    'JKHY',
    'RYAN',
    'PSKY',
    'CHWY',
    'UNM',
    'SNAP',
    'JEF',
    'OKLO',
    'HLI',
    'IVZ',
    'AGNC',
    'BMNR',
    'GLXY',
    'GMED',
    'DT',
    'ACGLO',
    'DOC',
    'SMMT',
    'EPAM',
    'JOBY',
    'CMA',
    'TEM',
    'NYT',
    'NTNX',
    'BSY',
    'DOCU',
    'BXP',
    'MDGL',
    'MICC',
    'QGEN',
    'WTRG',
    'SARO',
    'CRL',
    'UHAL',
    'MOH',
    'MANH',
    'FDS',
    'CART',
    'SEIC',
    'TECH',
    'CAE',
    'YMM',
    'PCOR',
    'KLAR',
    'ARMK',
    'JAZZ',
    'CHYM',
    'VNO-PL',
    'UHAL-B',
    'SAIL',
    'SANM',
    'BIO-B',
    'SOLS',
    'BROS',
    'REXR',
    'ABVX',
    'UWMC',
    'ARWR',
    'GS-PC',
    'RGEN',
    'DOX',
    'STEP',
    'MORN',
    'AMG',
    'UEC',
    'LUMN',
    'GGAL',
    'QBTS',
    'TTAN',
    'RZB',
    'AGNCM',
    'AGNCN',
    'PEGA',
    'UGI',
    'PL'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

# Remove now-unused trailing rows (previously up to row 428)
$ws.Range("A388:A428").ClearContents()

